$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginDetails")

# Duplicate row 2 (UserName/Password values) into a new row 3
$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2

# Select the new row, matching the entire-row selection left behind in Excel
$ws.Range("A3:XFD3").Select()
